# Apply the "Valid Species" list update:
#  - Insert "Pogoniulus atroflavus" in its correct alphabetical spot
#    (right before "Pogoniulus bilineatus").
#  - Remove the "Turdus philomelos" species group (5 rows):
#       Turdus philomelos
#       Turdus philomelos clarkei
#       Turdus philomelos hebridensis
#       Turdus philomelos nataliae
#       Turdus philomelos philomelos

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Valid Species")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# --- Find the row for "Pogoniulus bilineatus" so we can insert just above it ---
$insertRow = -1
for ($r = 1; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 1).Value() -eq "Pogoniulus bilineatus") {
        $insertRow = $r
        break
    }
}

if ($insertRow -ne -1) {
    $ws.Rows.Item($insertRow).Insert()
    $ws.Cells.Item($insertRow, 1).Value = "Pogoniulus atroflavus"
}

# --- Remove the "Turdus philomelos" group rows ---
$toRemove = @(
    "Turdus philomelos",
    "Turdus philomelos clarkei",
    "Turdus philomelos hebridensis",
    "Turdus philomelos nataliae",
    "Turdus philomelos philomelos"
)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = $lastRow; $r -ge 1; $r--) {
    $val = $ws.Cells.Item($r, 1).Value()
    if ($toRemove -contains $val) {
        $ws.Rows.Item($r).Delete()
    }
}
